$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 27 (ALMEIDA CUATIN JHONATHANN CARLOS / SANTANA JIMENEZ MARIA ELENA)
$ws1.Range("D27").Value = 7041.6
$ws1.Range("L27").Value = 2238.91

# Row 34 totals ("X de 32" counters)
$ws1.Range("D34").Value = "5 de 32"
$ws1.Range("L34").Value = "3 de 32"

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 27 (ALMEIDA CUATIN JHONATHANN CARLOS / SANTANA JIMENEZ MARIA ELENA), agosto
$ws2.Range("F27").Value = 11570.11

# Row 34 totals, agosto
$ws2.Range("F34").Value = 28259.35

# Column F width 13 -> 14 (ColumnWidth input tuned so OOXML width lands on 14)
$ws2.Columns.Item(6).ColumnWidth = 13.1666666666667

# ---------------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO
$ws3.Range("D3").Value = 14516.93
$ws3.Range("E3").Value = -11396.8155
$ws3.Range("F3").Value = 4.652691431676626

# Row 15: PIEDRA SINTERIZADA
$ws3.Range("D15").Value = 3445.75
$ws3.Range("E15").Value = -2918.72
$ws3.Range("F15").Value = 6.538052862265905

# Row 19: TOTAL
$ws3.Range("D19").Value = 28777.56
$ws3.Range("E19").Value = 3331.721075557873
$ws3.Range("F19").Value = 0.8962380668779897

# Column D width 13 -> 14
$ws3.Columns.Item(4).ColumnWidth = 13.1666666666667
